# Daily attendance processing - 2026-01-09 06:06:31
#
# For a specific set of rows in the "Recorded By" column (G), the two
# comma-separated author names were swapped in order (e.g.
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
# This script finds each affected row (by its current/expected value)
# and reverses the order of the two comma-separated entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows in column G whose "Recorded By" value needs its two
# comma-separated authors swapped.
$rowsToSwap = @(7,10,11,12,13,14,15,17,18,19,20,21,22,24,26,33,36,37,38,39,40,41,43,44,45,46,47,48,50,52,59,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rowsToSwap) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $current = $cell.Value2
    if ($current -ne $null -and $current -ne "") {
        $parts = $current -split ", ", 2
        if ($parts.Count -eq 2) {
            $cell.Value2 = $parts[1] + ", " + $parts[0]
        }
    }
}
